# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Source data: coinranking.com scrape -> Sheet1 columns B:Coin C:Link D:Price E:Volume(1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.767.04'
$ws.Range('E2').Value = '  +1.63%  '

$ws.Range('D3').Value = '3.093.72'

$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').Value = '''579.22'
$ws.Range('E5').Value = '  +1.53%  '

$ws.Range('D6').Value = '''172.36'
$ws.Range('E6').Value = '  +5.59%  '

$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('D8').Value = '3.088.97'
$ws.Range('E8').Value = '  +5.13%  '

$ws.Range('E9').Value = '  +1.36%  '

$ws.Range('D10').Value = '''6.43'
$ws.Range('E10').Value = '  -2.92%  '

$ws.Range('E11').Value = '  +3.08%  '

$ws.Range('E12').Value = '  +3.93%  '

$ws.Range('E13').Value = '  +2.07%  '

$ws.Range('D14').Value = '''37.23'

$ws.Range('E15').Value = '  +0.24%  '

$ws.Range('D16').Value = '3.604.53'
$ws.Range('E16').Value = '  +5.06%  '

$ws.Range('D17').Value = '66.738.25'
$ws.Range('E17').Value = '  +1.63%  '

$ws.Range('E18').Value = '  +1.72%  '

$ws.Range('D19').Value = '3.096.95'
$ws.Range('E19').Value = '  +5.26%  '

$ws.Range('D20').Value = '''16.33'
$ws.Range('E20').Value = '  +3.04%  '

$ws.Range('D21').Value = '''480.85'
$ws.Range('E21').Value = '  +7.97%  '

$ws.Range('E22').Value = '  +2.66%  '

$ws.Range('D23').Value = '''7.53'
$ws.Range('E23').Value = '  +3.69%  '

$ws.Range('E24').Value = '  +9.14%  '

$ws.Range('D25').Value = '''83.93'

$ws.Range('E26').Value = '  +5.49%  '

$ws.Range('D27').Value = '''10.04'
$ws.Range('E27').Value = '  +0.43%  '

$ws.Range('E28').Value = '  +0.05%  '

$ws.Range('E29').Value = '  -2.53%  '

$ws.Range('E30').Value = '  -1.54%  '

$ws.Range('E31').Value = '  +3.69%  '

# Row 32: coin identity swap (PEPE -> EthereumClassic)
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').Value = '''28.71'
$ws.Range('E32').Value = '  +5.30%  '

# Row 33: coin identity swap (EthereumClassic -> PEPE)
$ws.Range('B33').Value = 'PEPE'
$ws.Range('C33').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D33').Value = '''0.0000100'
$ws.Range('E33').Value = '  -0.17%  '

$ws.Range('E34').Value = '  -1.45%  '

$ws.Range('E35').Value = '  +0.07%  '

$ws.Range('E36').Value = '  +3.38%  '

$ws.Range('D37').Value = '''0.991'
$ws.Range('E37').Value = '  +2.08%  '

$ws.Range('D38').Value = '''48.25'
$ws.Range('E38').Value = '  +2.82%  '

$ws.Range('E39').Value = '  +7.01%  '

$ws.Range('E40').Value = '  +1.94%  '

$ws.Range('E41').Value = '  +4.98%  '

$ws.Range('E42').Value = '  +0.84%  '

$ws.Range('E43').Value = '  +2.20%  '

$ws.Range('D44').Value = '''2.79'
$ws.Range('E44').Value = '  -0.68%  '

$ws.Range('D45').Value = '2.839.88'
$ws.Range('E45').Value = '  +6.51%  '

$ws.Range('D46').Value = '''0.0359'
$ws.Range('E46').Value = '  +2.93%  '

$ws.Range('D47').Value = '''382.85'
$ws.Range('E47').Value = '  +0.57%  '

$ws.Range('D48').Value = '''135.21'
$ws.Range('E48').Value = '  +1.24%  '

$ws.Range('E49').Value = '  +0.02%  '

$ws.Range('D50').Value = '''24.87'
$ws.Range('E50').Value = '  +4.21%  '
